$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
$ws.Range("D3").Interior.Color = 255
"Set D3 interior"
$ws.Rows(3).RowHeight = 28.8
"Set row height"
